$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.757752656936646
$ws.Range("B1").Value = 2.063446521759033
$ws.Range("C1").Value = 2.189067363739014
$ws.Range("D1").Value = 2.409654140472412
$ws.Range("E1").Value = 3.05979585647583
